$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 1.83
$ws.Range("H3").Value = 3.4
$ws.Range("I3").Value = 4.75
$ws.Range("J3").Value = 1.13
$ws.Range("K3").Value = 6
$ws.Range("L3").Value = 1.62
$ws.Range("M3").Value = 2.3
$ws.Range("N3").Value = 2.88
$ws.Range("O3").Value = 1.4
$ws.Range("P3").Value = 1.62
$ws.Range("Q3").Value = 2.2
$ws.Range("R3").Value = 2.5
$ws.Range("S3").Value = 1.5
$ws.Range("T3").Value = 4.75
$ws.Range("U3").Value = 7
$ws.Range("W3").Value = 15
$ws.Range("X3").Value = 21
$ws.Range("Z3").Value = 6
$ws.Range("AE3").Value = 8.5
$ws.Range("AF3").Value = 21
$ws.Range("AG3").Value = 17
# Row 4
$ws.Range("G4").Value = 3
$ws.Range("I4").Value = 2.5
$ws.Range("J4").Value = 1.1
$ws.Range("K4").Value = 7
$ws.Range("L4").Value = 1.44
$ws.Range("M4").Value = 2.75
$ws.Range("N4").Value = 2.35
$ws.Range("O4").Value = 1.57
$ws.Range("P4").Value = 1.53
$ws.Range("Q4").Value = 2.38
$ws.Range("T4").Value = 8
$ws.Range("U4").Value = 13
$ws.Range("V4").Value = 12
$ws.Range("W4").Value = 34
$ws.Range("X4").Value = 29
$ws.Range("Z4").Value = 7
$ws.Range("AB4").Value = 15
$ws.Range("AC4").Value = 51
$ws.Range("AD4").Value = 401
$ws.Range("AE4").Value = 7
$ws.Range("AF4").Value = 11
$ws.Range("AG4").Value = 10
$ws.Range("AH4").Value = 26
$ws.Range("AI4").Value = 23
# Row 5
$ws.Range("G5").Value = 2.55
$ws.Range("H5").Value = 3.1
$ws.Range("I5").Value = 3
$ws.Range("J5").Value = 1.08
$ws.Range("K5").Value = 8
$ws.Range("N5").Value = 2.15
$ws.Range("O5").Value = 1.67
$ws.Range("T5").Value = 7.5
$ws.Range("U5").Value = 12
$ws.Range("V5").Value = 10
$ws.Range("W5").Value = 23
$ws.Range("X5").Value = 21
$ws.Range("Y5").Value = 34
$ws.Range("Z5").Value = 8
$ws.Range("AE5").Value = 8.5
$ws.Range("AF5").Value = 13
$ws.Range("AG5").Value = 11
$ws.Range("AH5").Value = 29
$ws.Range("AI5").Value = 26
$ws.Range("AJ5").Value = 34
# Row 6
$ws.Range("G6").Value = 2.25
$ws.Range("H6").Value = 3
$ws.Range("I6").Value = 3.6
$ws.Range("L6").Value = 1.53
$ws.Range("M6").Value = 2.5
$ws.Range("U6").Value = 9
$ws.Range("W6").Value = 21
$ws.Range("AB6").Value = 19
$ws.Range("AE6").Value = 8
$ws.Range("AF6").Value = 17
$ws.Range("AG6").Value = 13
$ws.Range("AI6").Value = 34
# Row 7
$ws.Range("H7").Value = 2.9
$ws.Range("I7").Value = 2.88
$ws.Range("J7").Value = 1.11
$ws.Range("K7").Value = 6.5
$ws.Range("L7").Value = 1.5
$ws.Range("M7").Value = 2.63
$ws.Range("N7").Value = 2.5
$ws.Range("O7").Value = 1.5
$ws.Range("P7").Value = 1.57
$ws.Range("Q7").Value = 2.25
$ws.Range("R7").Value = 2.05
$ws.Range("S7").Value = 1.7
$ws.Range("W7").Value = 29
$ws.Range("X7").Value = 26
$ws.Range("Z7").Value = 6.5
$ws.Range("AC7").Value = 67
$ws.Range("AD7").Value = 501
$ws.Range("AE7").Value = 7
$ws.Range("AF7").Value = 12
$ws.Range("AI7").Value = 29
# Row 8
$ws.Range("G8").Value = 1.22
$ws.Range("H8").Value = 6.25
$ws.Range("N8").Value = 1.65
$ws.Range("O8").Value = 2.2
$ws.Range("P8").Value = 1.3
$ws.Range("Q8").Value = 3.4
$ws.Range("R8").Value = 2.38
$ws.Range("S8").Value = 1.53
$ws.Range("U8").Value = 5.5
$ws.Range("V8").Value = 10
$ws.Range("W8").Value = 6.5
$ws.Range("X8").Value = 12
$ws.Range("Z8").Value = 13
$ws.Range("AC8").Value = 101
$ws.Range("AE8").Value = 26
$ws.Range("AF8").Value = 51
$ws.Range("AG8").Value = 34
# Row 10
$ws.Range("N10").Value = 2.6
$ws.Range("O10").Value = 1.48
$ws.Range("R10").Value = 2.1
$ws.Range("S10").Value = 1.67
$ws.Range("T10").Value = 6
# Row 11
$ws.Range("H11").Value = 3
$ws.Range("J11").Value = 1.1
$ws.Range("L11").Value = 1.5
$ws.Range("M11").Value = 2.37
$ws.Range("N11").Value = 2.7
$ws.Range("O11").Value = 1.44
$ws.Range("P11").Value = 1.62
$ws.Range("Q11").Value = 2.2
$ws.Range("R11").Value = 2.2
$ws.Range("S11").Value = 1.62
$ws.Range("W11").Value = 19
$ws.Range("X11").Value = 21
$ws.Range("Z11").Value = 6
# Row 12
$ws.Range("J12").Value = 1.04
$ws.Range("K12").Value = 10
$ws.Range("L12").Value = 1.33
$ws.Range("V12").Value = 9
$ws.Range("Y12").Value = 34
$ws.Range("Z12").Value = 8.5
# Row 13
$ws.Range("G13").Value = 2.7
$ws.Range("H13").Value = 2.75
$ws.Range("I13").Value = 3
$ws.Range("J13").Value = 1.17
$ws.Range("K13").Value = 5
$ws.Range("T13").Value = 6
$ws.Range("U13").Value = 11
$ws.Range("V13").Value = 12
$ws.Range("W13").Value = 29
$ws.Range("Y13").Value = 51
$ws.Range("AE13").Value = 6.5
$ws.Range("AF13").Value = 13
# Row 20
$ws.Range("J20").Value = 1.11
$ws.Range("K20").Value = 6.5
$ws.Range("V20").Value = 9.5
$ws.Range("AC20").Value = 101
# Row 21
$ws.Range("G21").Value = 1.22
$ws.Range("H21").Value = 5.5
$ws.Range("J21").Value = 1.04
$ws.Range("K21").Value = 13
$ws.Range("L21").Value = 1.22
$ws.Range("M21").Value = 4
$ws.Range("N21").Value = 1.73
$ws.Range("O21").Value = 2.08
$ws.Range("P21").Value = 1.33
$ws.Range("Q21").Value = 3.25
$ws.Range("T21").Value = 6
$ws.Range("W21").Value = 7
$ws.Range("Y21").Value = 41
$ws.Range("Z21").Value = 11
$ws.Range("AA21").Value = 11
$ws.Range("AE21").Value = 26
$ws.Range("AF21").Value = 51
$ws.Range("AG21").Value = 34
$ws.Range("AH21").Value = 201
$ws.Range("AI21").Value = 101
$ws.Range("AJ21").Value = 81
# Row 22
$ws.Range("G22").Value = 1.95
$ws.Range("I22").Value = 3.9
$ws.Range("U22").Value = 8
$ws.Range("AE22").Value = 9
$ws.Range("AF22").Value = 19
$ws.Range("AG22").Value = 15
$ws.Range("AI22").Value = 41
# Row 23
$ws.Range("G23").Value = 1.93
$ws.Range("H23").Value = 2.9
$ws.Range("I23").Value = 4.4
$ws.Range("Q23").Value = 2.1
$ws.Range("R23").Value = 2.05
$ws.Range("T23").Value = 5.3
$ws.Range("U23").Value = 8
$ws.Range("W23").Value = 17
$ws.Range("X23").Value = 19
$ws.Range("Z23").Value = 6.1
$ws.Range("AA23").Value = 5.9
$ws.Range("AB23").Value = 18.5
$ws.Range("AE23").Value = 9
$ws.Range("AF23").Value = 23
$ws.Range("AH23").Value = 80
$ws.Range("AI23").Value = 55
# Row 24
$ws.Range("I24").Value = 3.2
$ws.Range("N24").Value = 2.05
$ws.Range("O24").Value = 1.75
$ws.Range("X24").Value = 19
# Row 25
$ws.Range("G25").Value = 3.3
$ws.Range("H25").Value = 3.4
$ws.Range("I25").Value = 2
$ws.Range("J25").Value = 1.06
$ws.Range("K25").Value = 10
$ws.Range("T25").Value = 9.5
$ws.Range("U25").Value = 17
$ws.Range("X25").Value = 29
$ws.Range("AA25").Value = 6.5
$ws.Range("AB25").Value = 15
# Row 41
$ws.Range("G41").Value = 2.88
$ws.Range("I41").Value = 2.6
$ws.Range("V41").Value = 11
$ws.Range("W41").Value = 29
$ws.Range("AI41").Value = 19
$ws.Range("AJ41").Value = 26
# Row 42
$ws.Range("L42").Value = 1.33
$ws.Range("M42").Value = 3.25
# Row 43
$ws.Range("H43").Value = 4.1
$ws.Range("I43").Value = 1.65
$ws.Range("K43").Value = 21
$ws.Range("P43").Value = 1.22
$ws.Range("Q43").Value = 4
$ws.Range("V43").Value = 17
# Row 46
$ws.Range("G46").Value = 4.2
$ws.Range("H46").Value = 3.25
$ws.Range("I46").Value = 1.91
$ws.Range("N46").Value = 2.25
$ws.Range("O46").Value = 1.62
$ws.Range("R46").Value = 2
$ws.Range("S46").Value = 1.73
$ws.Range("T46").Value = 10
$ws.Range("U46").Value = 21
$ws.Range("V46").Value = 15
$ws.Range("W46").Value = 41
$ws.Range("Y46").Value = 41
$ws.Range("AB46").Value = 17
$ws.Range("AE46").Value = 6
$ws.Range("AF46").Value = 8
$ws.Range("AH46").Value = 15
